# Rebuild the worksheet with the new "invalid_template" column layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, extend the bold/bordered header style (already applied to A1:F1)
# to the new header columns G1:J1 *before* we touch any cell content, so
# the original style index (s="1") is reused rather than a brand-new
# (and non-identical) style being registered.
$ws.Range("A1").Copy()
$ws.Range("G1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Helper: writes a text value into a cell while preventing Excel's
# automatic number/date inference (e.g. "2024-01-01" -> date serial,
# "11050501" -> numeric). We compute the literal text with a TEXT()
# formula in a scratch cell far outside the used range, then copy only
# the *value* (PasteSpecial values) into the destination cell. This
# keeps the value a plain string without registering a new cell style.
function Set-TextValue {
    param($cellRef, [string]$text)
    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '=TEXT("' + $escaped + '","@")'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "document_id"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "account_code"
$ws.Range("D1").Value = "movement"
$ws.Range("E1").Value = "customer_identification"
$ws.Range("F1").Value = "branch_office"
$ws.Range("G1").Value = "description"
$ws.Range("H1").Value = "cost_center"
$ws.Range("I1").Value = "value"
$ws.Range("J1").Value = "observations"

# ---- Row 2 ----
$ws.Range("A2").Value = 27441
Set-TextValue "B2" "2024-01-01"
Set-TextValue "C2" "11050501"
$ws.Range("D2").Value = "Debit"
Set-TextValue "E2" "13832081"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "Test debit"
$ws.Range("H2").Value = 235
$ws.Range("I2").Value = 119000
$ws.Range("J2").Value = "Observaciones"

# ---- Row 3 ----
$ws.Range("A3").Value = 27441
Set-TextValue "B3" "2024-01-01"
Set-TextValue "C3" "11100501"
$ws.Range("D3").Value = "Credit"
Set-TextValue "E3" "13832081"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = "Test credit"
$ws.Range("H3").Value = 235
$ws.Range("I3").Value = 90000
$ws.Range("J3").Value = "Observaciones"
